$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected "Fitness" (column C) values for rows 2 through 128, reflecting
# the fix to the SA algorithm referenced in the commit message
# ("correction in sa algorithm and 746 logs"). All other rows/columns in
# the log are unchanged.
$newValues = @(
    8320,8320,8320,8320,8295,8295,8295,8295,8295,8295,8295,8295,8295,8295,8295,
    8295,8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,
    8028,8028,8028,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
